# Actualización automática 2025-09-11 17:30:08
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("I16").Value = 25.2
$wsGrupo.Range("D50").Value = 183.16
$wsGrupo.Range("I54").Value = "10 de 52"

# --- Sheet "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F16").Value = 25.2
$wsMensual.Range("F51").Value = 575.5599999999999
$wsMensual.Range("F52").Value = 575.5599999999999
$wsMensual.Range("F58").Value = 29079.49

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D3").Value = 1764.28
$wsCumpl.Range("E3").Value = 15904.8670988183
$wsCumpl.Range("F3").Value = 0.09985088641420581

$wsCumpl.Range("D7").Value = 2349.9
$wsCumpl.Range("E7").Value = -1463.188983712426
$wsCumpl.Range("F7").Value = 2.650130602683176

$wsCumpl.Range("D15").Value = 28075.1
$wsCumpl.Range("E15").Value = 93979.73551083436
$wsCumpl.Range("F15").Value = 0.2300203829082042
